# Auto-generated Excel COM-interop script applying the scheduled market-data refresh
# to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 231.27272
$ws.Range("I5").Value = 231.27272
$ws.Range("K5").Value = 231.27272
$ws.Range("M5").Value = -116.27272
$ws.Range("H6").Value = 97.5
$ws.Range("I6").Value = 97.5
$ws.Range("K6").Value = 292.5
$ws.Range("M6").Value = -180.5
$ws.Range("H33").Value = 450.1875
$ws.Range("I33").Value = 450.1875
$ws.Range("K33").Value = 450.1875
$ws.Range("M33").Value = -221.1875
$ws.Range("H55").Value = 359.47058
$ws.Range("I55").Value = 281.3
$ws.Range("J55").Value = 471.14285
$ws.Range("K55").Value = 281.3
$ws.Range("L55").Value = 471.14285
$ws.Range("M55").Value = -67.30000000000001
$ws.Range("N55").Value = -899.14285
$ws.Range("H80").Value = 1560.7
$ws.Range("I80").Value = 881.6667
$ws.Range("J80").Value = 1851.7142
$ws.Range("K80").Value = 2645.0001
$ws.Range("L80").Value = 5555.142599999999
$ws.Range("M80").Value = -1647.0001
$ws.Range("N80").Value = -7551.142599999999
$ws.Range("H83").Value = 1560.7
$ws.Range("I83").Value = 881.6667
$ws.Range("J83").Value = 1851.7142
$ws.Range("K83").Value = 7935.0003
$ws.Range("L83").Value = 16665.4278
$ws.Range("M83").Value = -2943.0003
$ws.Range("N83").Value = -26649.4278
$ws.Range("H111").Value = 10432.333
$ws.Range("I111").Value = 10432.333
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 31296.999
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -28229.999
$ws.Range("H116").Value = 5320.2
$ws.Range("I116").Value = 4585.5713
$ws.Range("J116").Value = 5963
$ws.Range("K116").Value = 4585.5713
$ws.Range("L116").Value = 5963
$ws.Range("M116").Value = -1143.5713
$ws.Range("N116").Value = -12847
$ws.Range("H137").Value = 1165624.9
$ws.Range("I137").Value = 1517507
$ws.Range("J137").Value = 4413.8
$ws.Range("K137").Value = 4552521
$ws.Range("L137").Value = 13241.4
$ws.Range("M137").Value = -4549971
$ws.Range("N137").Value = -18341.4
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 30048
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("H63").Value = 2165.5715
$ws.Range("I63").Value = 2068.1667
$ws.Range("J63").Value = 2750
$ws.Range("K63").Value = 2068.1667
$ws.Range("L63").Value = 2750
$ws.Range("M63").Value = -1382.1667
$ws.Range("N63").Value = -4122
$ws.Range("H66").Value = 2165.5715
$ws.Range("I66").Value = 2068.1667
$ws.Range("J66").Value = 2750
$ws.Range("K66").Value = 10340.8335
$ws.Range("L66").Value = 13750
$ws.Range("M66").Value = -6908.833500000001
$ws.Range("N66").Value = -20614
$ws.Range("H74").Value = 348689.5
$ws.Range("J74").Value = 2057.125
$ws.Range("L74").Value = 2057.125
$ws.Range("N74").Value = -3805.125
$ws.Range("H77").Value = 348689.5
$ws.Range("J77").Value = 2057.125
$ws.Range("L77").Value = 10285.625
$ws.Range("N77").Value = -19021.625
$ws.Range("H80").Value = 72477.5
$ws.Range("J80").Value = 89970
$ws.Range("L80").Value = 89970
$ws.Range("N80").Value = -91966
$ws.Range("H83").Value = 72477.5
$ws.Range("J83").Value = 89970
$ws.Range("L83").Value = 269910
$ws.Range("N83").Value = -279894
$ws.Range("H122").Value = 3960.5386
$ws.Range("I122").Value = 4250
$ws.Range("J122").Value = 3873.7
$ws.Range("K122").Value = 12750
$ws.Range("L122").Value = 11621.1
$ws.Range("M122").Value = -10300
$ws.Range("N122").Value = -16521.1
$ws.Range("H132").Value = 2009.2727
$ws.Range("I132").Value = 1805
$ws.Range("K132").Value = 5415
$ws.Range("M132").Value = -2885
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 926.875
$ws.Range("J64").Value = 1424
$ws.Range("L64").Value = 1424
$ws.Range("N64").Value = -1874
$ws.Range("H67").Value = 926.875
$ws.Range("J67").Value = 1424
$ws.Range("L67").Value = 1424
$ws.Range("N67").Value = -2984
$ws.Range("H105").Value = 20001996
$ws.Range("J105").Value = 41669316
$ws.Range("L105").Value = 41669316
$ws.Range("N105").Value = -41672810
$ws.Range("H132").Value = 105298.664
$ws.Range("J132").Value = 105298.664
$ws.Range("L132").Value = 105298.664
$ws.Range("N132").Value = -115418.664
$ws.Range("H134").Value = 11897.111
$ws.Range("I134").Value = 4343.7
$ws.Range("J134").Value = 21338.875
$ws.Range("K134").Value = 13031.1
$ws.Range("L134").Value = 64016.625
$ws.Range("M134").Value = -10496.1
$ws.Range("N134").Value = -69086.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3913.5178
$ws.Range("I31").Value = 3317
$ws.Range("K31").Value = 3317
$ws.Range("M31").Value = -3022
$ws.Range("H34").Value = 3913.5178
$ws.Range("I34").Value = 3317
$ws.Range("K34").Value = 3317
$ws.Range("M34").Value = -3115
$ws.Range("H58").Value = 1994.973
$ws.Range("I58").Value = 1475.4736
$ws.Range("J58").Value = 2543.3333
$ws.Range("K58").Value = 1475.4736
$ws.Range("L58").Value = 2543.3333
$ws.Range("M58").Value = -1272.4736
$ws.Range("N58").Value = -2949.3333
$ws.Range("H62").Value = 12503649
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 12503649
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H99").Value = 4257.75
$ws.Range("I99").Value = 4614
$ws.Range("J99").Value = 1764
$ws.Range("K99").Value = 4614
$ws.Range("L99").Value = 1764
$ws.Range("M99").Value = -3116
$ws.Range("N99").Value = -4760
$ws.Range("H126").Value = 4257.75
$ws.Range("I126").Value = 4614
$ws.Range("J126").Value = 1764
$ws.Range("K126").Value = 13842
$ws.Range("L126").Value = 5292
$ws.Range("M126").Value = -11372
$ws.Range("N126").Value = -10232
$ws.Range("H132").Value = 5380679
$ws.Range("I132").Value = 4390.383
$ws.Range("J132").Value = 22226384
$ws.Range("K132").Value = 13171.149
$ws.Range("L132").Value = 66679152
$ws.Range("M132").Value = -10641.149
$ws.Range("N132").Value = -66684212
$ws.Range("H136").Value = 1994.973
$ws.Range("I136").Value = 1475.4736
$ws.Range("J136").Value = 2543.3333
$ws.Range("K136").Value = 4426.4208
$ws.Range("L136").Value = 7629.999899999999
$ws.Range("M136").Value = -1876.4208
$ws.Range("N136").Value = -12729.9999
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 151.6923
$ws.Range("I12").Value = 32.666668
$ws.Range("K12").Value = 98.000004
$ws.Range("M12").Value = 74.999996
$ws.Range("H28").Value = 900
$ws.Range("I28").Value = 900
$ws.Range("K28").Value = 2700
$ws.Range("M28").Value = -2468

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2891.4285
$ws.Range("I31").Value = 873.3333
$ws.Range("J31").Value = 15000
$ws.Range("K31").Value = 873.3333
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = -581.3333
$ws.Range("N31").Value = -15584
$ws.Range("H37").Value = 2891.4285
$ws.Range("I37").Value = 873.3333
$ws.Range("J37").Value = 15000
$ws.Range("K37").Value = 873.3333
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = -596.3333
$ws.Range("N37").Value = -15554
$ws.Range("H113").Value = 8438.429
$ws.Range("I113").Value = 3899.8
$ws.Range("K113").Value = 3899.8
$ws.Range("M113").Value = -1729.8
$ws.Range("H132").Value = 2430.1667
$ws.Range("I132").Value = 2358.9375
$ws.Range("K132").Value = 7076.8125
$ws.Range("M132").Value = -4546.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3561.75
$ws.Range("I22").Value = 2717.625
$ws.Range("J22").Value = 5250
$ws.Range("K22").Value = 2717.625
$ws.Range("L22").Value = 5250
$ws.Range("M22").Value = -2422.625
$ws.Range("N22").Value = -5840
$ws.Range("H27").Value = 3561.75
$ws.Range("I27").Value = 2717.625
$ws.Range("J27").Value = 5250
$ws.Range("K27").Value = 2717.625
$ws.Range("L27").Value = 5250
$ws.Range("M27").Value = -2610.625
$ws.Range("N27").Value = -5464
$ws.Range("H40").Value = 87035.914
$ws.Range("I40").Value = 94493.73
$ws.Range("K40").Value = 94493.73
$ws.Range("M40").Value = -94357.73
$ws.Range("H68").Value = 4638.8
$ws.Range("I68").Value = 3298.75
$ws.Range("K68").Value = 3298.75
$ws.Range("M68").Value = -2549.75
$ws.Range("H71").Value = 4638.8
$ws.Range("I71").Value = 3298.75
$ws.Range("K71").Value = 16493.75
$ws.Range("M71").Value = -12749.75
$ws.Range("H132").Value = 6230.7144
$ws.Range("I132").Value = 6428.143
$ws.Range("K132").Value = 19284.429
$ws.Range("M132").Value = -16754.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 19006
$ws.Range("I26").Value = 19006
$ws.Range("K26").Value = 19006
$ws.Range("M26").Value = -18713
$ws.Range("H62").Value = 6070.5454
$ws.Range("J62").Value = 5775.857
$ws.Range("L62").Value = 5775.857
$ws.Range("N62").Value = -7023.857
$ws.Range("H65").Value = 6070.5454
$ws.Range("J65").Value = 5775.857
$ws.Range("L65").Value = 28879.285
$ws.Range("N65").Value = -35119.285
$ws.Range("H92").Value = 113249.5
$ws.Range("J92").Value = 113249.5
$ws.Range("L92").Value = 113249.5
$ws.Range("N92").Value = -118241.5
$ws.Range("H132").Value = 6647
$ws.Range("I132").Value = 6950.4614
$ws.Range("K132").Value = 20851.3842
$ws.Range("M132").Value = -18321.3842
